$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Attack Up" / "Attack" bonus block that used to live in rows 8:13 is being
# removed. The rows below it (the "Speed Up", "MaxHP Up" and "Combo Up" blocks,
# previously rows 14:27) move up to take its place. Row (ID) numbers in column A
# are untouched - only columns B:E carry data that needs to shift.
#
# Target content for rows 8:21 (after the shift):
$rows = @(
    # row, BonusName,  BonusValue, Rarelity, RarelityStyled, Type,    TypeStyled
    @{R=8;  B="Speed Up";  C=30;  D="Normal"; Ds=$true;  E="Speed"; Es=$true},
    @{R=9;  B="Speed Up";  C=60;  D="Normal"; Ds=$true;  E="Speed"; Es=$true},
    @{R=10; B="Speed Up";  C=90;  D="Rare";   Ds=$true;  E="Speed"; Es=$true},
    @{R=11; B="Speed Up";  C=120; D="Rare";   Ds=$true;  E="Speed"; Es=$true},
    @{R=12; B="Speed Up";  C=150; D="Unique"; Ds=$true;  E="Speed"; Es=$true},
    @{R=13; B="Speed Up";  C=300; D="Legend"; Ds=$true;  E="Speed"; Es=$true},
    @{R=14; B="MaxHP Up";  C=10;  D="Normal"; Ds=$true;  E="MaxHP"; Es=$false},
    @{R=15; B="MaxHP Up";  C=20;  D="Normal"; Ds=$true;  E="MaxHP"; Es=$false},
    @{R=16; B="MaxHP Up";  C=30;  D="Rare";   Ds=$true;  E="MaxHP"; Es=$false},
    @{R=17; B="MaxHP Up";  C=40;  D="Rare";   Ds=$true;  E="MaxHP"; Es=$false},
    @{R=18; B="MaxHP Up";  C=50;  D="Unique"; Ds=$true;  E="MaxHP"; Es=$false},
    @{R=19; B="MaxHP Up";  C=100; D="Legend"; Ds=$false; E="MaxHP"; Es=$false},
    @{R=20; B="Combo Up";  C=1;   D="Rare";   Ds=$false; E="Combo"; Es=$false},
    @{R=21; B="Combo Up";  C=2;   D="Unique"; Ds=$true;  E="Combo"; Es=$false}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    if (-not $row.Ds) { $ws.Range("D$r").ClearFormats() | Out-Null }
    if (-not $row.Es) { $ws.Range("E$r").ClearFormats() | Out-Null }
}

# Rows 22:27 are now vacated entirely (no more bonus entries below Combo Up) -
# B:E are fully cleared (content + formatting) while the ID cells in column A
# keep their (now empty) formatting.
$ws.Range("B22:E27").Clear() | Out-Null
$ws.Range("A22:A27").ClearContents() | Out-Null

$ws.Range("E16").Select() | Out-Null
